$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '4.69%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-4.50%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.212'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.32%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05885'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.94%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.704'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.53%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8708'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.15%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9774'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '13.49%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01050'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1,656.88%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1412'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.86%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07186'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.58%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03172'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.12%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09221'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.71%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001545'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.22%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006040'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.76%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.499'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.26%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.219'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.37%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.205'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.07%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3174'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.07%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03491'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.87%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1290'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.58%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.528'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.29%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04169'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.37%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1364'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.69%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001222'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.11%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004564'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '10.05%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001198'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.12%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001470'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.51%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03833'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.38%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005473'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-4.14%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1104'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.56%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002341'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.39%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009542'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '1.44%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005408'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.44%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.11%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.09488'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '5.57%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002128'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-13.23%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002097'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.11%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001997'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.11%'
